# Update countries & provincias Spain
# Applies the refreshed COVID-19 snapshot data to the "Pais" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Datos actualizados a ..." footer timestamp.
$ws.Range("A1").Value = "Datos actualizados a 9 de Octubre de 2020 a las 01:35"

# 2. Two row pairs swapped which country name they show (the underlying
#    source re-ordered these two pairs while all other rows kept their
#    country label). Row numbers / data otherwise stay fixed.
$ws.Range("A45").Value = "Chequia"
$ws.Range("A46").Value = "Nepal"

$ws.Range("A95").Value = "Noruega"
$ws.Range("A96").Value = "Senegal"

# 3. Refresh the numeric columns (B=Casos totales, C=Nuevos casos,
#    D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy,
#    H=Muertes) for every row whose figures changed.

function Set-Row($row, $b, $c, $d, $e, $g, $h) {
    if ($null -ne $b) { $ws.Cells.Item($row, 2).Value = $b }
    if ($null -ne $c) { $ws.Cells.Item($row, 3).Value = $c }
    if ($null -ne $d) { $ws.Cells.Item($row, 4).Value = $d }
    if ($null -ne $e) { $ws.Cells.Item($row, 5).Value = $e }
    if ($null -ne $g) { $ws.Cells.Item($row, 7).Value = $g }
    if ($null -ne $h) { $ws.Cells.Item($row, 8).Value = $h }
}

# Row 4  - Estados Unidos
Set-Row 4 7829136 52273 5015621 2595877 857 217638

# Row 6  - Brasil
Set-Row 6 5029539 27182 $null 465941 730 149034

# Row 10 - Argentina
Set-Row 10 856369 15454 684844 148815 484 22710

# Row 11 - Peru
Set-Row 11 838614 2952 728216 77300 89 33098

# Row 37 - Panama
Set-Row 37 118054 754 94391 21200 15 2463

# Row 45 - Nepal (label) -> now shows Chequia's refreshed figures
Set-Row 45 100757 5397 51160 48728 40 869

# Row 46 - Chequia (label) -> now shows Nepal's previous figures
Set-Row 46 98617 4364 71343 26684 12 590

# Row 60 - Nigeria
Set-Row 60 59841 103 51551 7177 $null $null

# Row 81 - Australia
Set-Row 81 27206 24 24951 1358 $null $null

# Row 84 - Bulgaria
Set-Row 84 23259 516 15563 6816 7 880

# Row 88 - Camerun
Set-Row 88 21203 279 20117 663 3 423

# Row 95 - Senegal (label) -> now shows Noruega's refreshed figures
Set-Row 95 15221 209 11863 3083 $null 275

# Row 96 - Noruega (label) -> now shows Senegal's previous figures
Set-Row 96 15190 16 13068 1809 $null 313

# Row 99 - Sudan
Set-Row 99 13670 2 $null 6070 $null $null

# Row 103 - Guinea
Set-Row 103 10901 38 10232 601 $null $null

# Row 104 - Consejo Danes para los Refugiados
Set-Row 104 10822 18 10242 304 $null $null

# Row 112 - Haiti
Set-Row 112 8854 16 $null 1611 1 230

# Row 115 - Mauritania
Set-Row 115 7540 5 7227 151 $null $null

# Row 131 - Trinidad yTobago
Set-Row 131 4963 76 3076 1801 2 86

# Row 158 - Uruguay
Set-Row 158 2226 20 1904 273 $null $null
